$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.249.53"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "'3.796.41"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'668.81"
$ws.Range("E5").Value = "  +6.78%  "
$ws.Range("D6").Value = "'168.48"
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("D7").Value = "'3.795.16"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.462"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'7.06"
$ws.Range("E12").Value = "  +4.67%  "
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "'35.59"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "'4.437.26"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "'3.798.81"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "'70.267.34"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "'17.63"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "'11.40"
$ws.Range("E21").Value = "  +18.56%  "
$ws.Range("D22").Value = "'474.50"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").Value = "'0.712"
$ws.Range("D24").Value = "'83.26"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  -4.97%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").Value = "'10.27"
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("D28").Value = "'2.11"
$ws.Range("E28").Value = "  -2.87%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'3.947.99"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("E31").Value = "  +5.69%  "
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("D33").Value = "'7.40"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("D34").Value = "'29.50"
$ws.Range("E34").Value = "  +2.05%  "
$ws.Range("D35").Value = "'0.177"
$ws.Range("E35").Value = "  +7.78%  "
$ws.Range("D36").Value = "'9.10"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'3.753.96"
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "'3.37"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").Value = "'5.94"
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.960"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "'2.11"
$ws.Range("E44").Value = "  +10.05%  "
$ws.Range("D46").Value = "'45.43"
$ws.Range("E46").Value = "  +5.14%  "
$ws.Range("D47").Value = "'159.43"
$ws.Range("E47").Value = "  +4.20%  "
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").Value = "'1.41"
$ws.Range("E50").Value = "  +4.21%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000288"
$ws.Range("E51").Value = "  +3.93%  "
